$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new "extraction line" rows (line7, line8) are inserted right after
# line6 / before extr1, which pushes the eight existing "extr#" rows down
# by two rows (old rows 8-15 -> new rows 10-17).
$ws.Rows("8:9").Insert()

# The insert carries column A's bold/bordered/centered style onto the new
# rows 8-9 automatically, but (in this runtime) also leaves a stray unused
# style definition behind. Re-apply the formatting explicitly from an
# existing cell so the new rows end up on the very same style index that
# the rest of column A already uses.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the two new rows: line7, line8.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Fix up the sequential index column for the rows that were pushed down
# (it's a running counter, not part of the shifted row content).
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# in_service flips to TRUE for what are now extr2 (row 11) and extr3 (row 12).
$ws.Range("E11").Value = $true
$ws.Range("E12").Value = $true
